# Insert a new data row at row 204 (pushing the existing rows 204-304 down to 205-305)
# and populate the new row with the weekly record added by this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(204).Insert()

$ws.Cells.Item(204, 1).Value = 6
$ws.Cells.Item(204, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(204, 3).Value = "Metropolitana"
$ws.Cells.Item(204, 4).Value = 45029
$ws.Cells.Item(204, 5).Value = 13
$ws.Cells.Item(204, 6).Value = 100112001
$ws.Cells.Item(204, 7).Value = "Berenjena"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 580
$ws.Cells.Item(204, 11).Value = 6000
$ws.Cells.Item(204, 12).Value = 7000
$ws.Cells.Item(204, 13).Value = 6448
$ws.Cells.Item(204, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(204, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(204, 16).Value = 107
$ws.Cells.Item(204, 17).Value = 60
$ws.Cells.Item(204, 18).Value = "Hortaliza"
